# Split two runs in the "Свойства" document so that the level numbers
# 7 and 8 (and the word "фиолетовым") are rendered in purple (7030A0),
# matching the rest of the text's normal color.
#
# Strategy: locate the whole original (unsplit) run text with Find, then
# carve out the exact character sub-ranges that must turn purple and set
# their Font.Color. Word automatically splits the underlying <w:r> run
# whenever a sub-range's resolved formatting differs from its neighbours,
# so no manual run/XML surgery is required.

$d = $word.ActiveDocument

# RGB(0x70, 0x30, 0xA0) expressed the way Word's object model expects
# (R + G*256 + B*65536) -> produces w:color val="7030A0" in the OOXML.
$purple = 10498160

# ---------------------------------------------------------------------
# 1) "Первые 6 уровней свойств могут быть накрафчены, а вот 7 и 8 могут
#    только выпасть. У предмета 5 уровня, может быть 7 уровень свойства
#    и у предмета 6 уровня может быть свойство 7 или 8 уровня"
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$found1 = $rng1.Find.Execute( `
    "Первые 6 уровней свойств могут быть накрафчены, а вот 7 и 8 могут только выпасть. У предмета 5 уровня, может быть 7 уровень свойства и у предмета 6 уровня может быть свойство 7 или 8 уровня", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $base1 = $rng1.Start

    # " 7 " right after "...а вот"
    $d.Range($base1 + 53, $base1 + 56).Font.Color = $purple
    # " 8 " right after the lone "и"
    $d.Range($base1 + 57, $base1 + 60).Font.Color = $purple
    # the lone "7" in "...может быть свойство 7 или 8 уровня"
    $d.Range($base1 + 175, $base1 + 176).Font.Color = $purple
    # the lone "8" in "...7 или 8 уровня"
    $d.Range($base1 + 181, $base1 + 182).Font.Color = $purple
}

# ---------------------------------------------------------------------
# 2) ". Такие свойства выделяются фиолетовым текстом."
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute( `
    ". Такие свойства выделяются фиолетовым текстом.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $base2 = $rng2.Start

    # "фиолетовым " (including the trailing space, excluding "текстом.")
    $d.Range($base2 + 28, $base2 + 39).Font.Color = $purple
}
